$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows (line7, line8) are inserted right after the "line6" row,
# shifting the former "extr1".."extr8" rows down by two. We recreate this
# by first copying the formatting (column A style) of the last existing
# data row down into the two brand-new rows 16 and 17, then writing the
# complete final table of values for every data row (2-17).

$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# name | A | B | C | D | E
$data = @(
  @(2,  0, "line1", 7,  9,  $true),
  @(3,  1, "line2", 9,  8,  $true),
  @(4,  2, "line3", 8,  10, $false),
  @(5,  3, "line4", 8,  11, $true),
  @(6,  4, "line5", 10, 5,  $true),
  @(7,  5, "line6", 12, 8,  $true),
  @(8,  6, "line7", 14, 11, $true),
  @(9,  7, "line8", 16, 9,  $true),
  @(10, 8, "extr1", 5,  12, $true),
  @(11, 9, "extr2", 5,  9,  $true),
  @(12, 10, "extr3", 10, 11, $false),
  @(13, 11, "extr4", 7,  8,  $false),
  @(14, 12, "extr5", 9,  11, $true),
  @(15, 13, "extr6", 7,  11, $false),
  @(16, 14, "extr7", 5,  7,  $false),
  @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}
